$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2 value
$ws.Range("G2").Value = 1339.26

# Add new row 3, copying row 2 content with modifications
$ws.Range("A3").Value = 342456
$ws.Range("B3").Value = "Thiago Aparecido Nogueira Basso"
$ws.Range("C3").Value = "Novembro2025"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "00001"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "r"
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 232.2
$ws.Range("H3").Value = "406 DE 03/11/2025 PROGESP"
$ws.Range("I3").Value = "407 DE 03/11/2025 PROGESP"

# Copy style from G2 to G3 (currency format)
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122) # xlPasteFormats
